$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name reflects "Through" date)
$ws.Name = "Through 2021-10-28"

# Update row 12 ("October (through 10-27/28)")
$ws.Range("A12").Value = "October (through 10-28)"
$ws.Range("C12").Value = 26
$ws.Range("D12").Value = 0.07140000000000001
$ws.Range("F12").Value = 42
$ws.Range("G12").Value = 0.1064
$ws.Range("H12").Value = 11
$ws.Range("I12").Value = 64
$ws.Range("J12").Value = 0.1467
$ws.Range("L12").Value = 54
$ws.Range("M12").Value = 0.0847
$ws.Range("O12").Value = 52
$ws.Range("P12").Value = 0.07140000000000001
$ws.Range("R12").Value = 135
$ws.Range("S12").Value = 0.0074
$ws.Range("U12").Value = 175

# Update row 13 ("Total")
$ws.Range("C13").Value = 222
$ws.Range("D13").Value = 0.126
$ws.Range("F13").Value = 425
$ws.Range("G13").Value = 0.1071
$ws.Range("H13").Value = 61
$ws.Range("I13").Value = 641
$ws.Range("J13").Value = 0.08690000000000001
$ws.Range("L13").Value = 541
$ws.Range("M13").Value = 0.1087
$ws.Range("O13").Value = 431
$ws.Range("P13").Value = 0.0983
$ws.Range("R13").Value = 983
$ws.Range("S13").Value = 0.0521
$ws.Range("U13").Value = 1340
